$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date-column style (s="1") down to the new rows (62:80) by copying
# the formatting from the last originally-styled date cell (A61).
$ws.Range("A61").Copy($ws.Range("A62:A80"))

$dates = @(
  "2017-10",
  "2017-11",
  "2017-12",
  "2017-01",
  "2017-02",
  "2017-03",
  "2017-04",
  "2017-05",
  "2017-06",
  "2017-07",
  "2017-08",
  "2017-09",
  "2018-10",
  "2018-11",
  "2018-12",
  "2018-01",
  "2018-02",
  "2018-03",
  "2018-04",
  "2018-05",
  "2018-06",
  "2018-07",
  "2018-08",
  "2018-09",
  "2019-10",
  "2019-11",
  "2019-12",
  "2019-01",
  "2019-02",
  "2019-03",
  "2019-04",
  "2019-05",
  "2019-06",
  "2019-07",
  "2019-08",
  "2019-09",
  "2020-10",
  "2020-11",
  "2020-12",
  "2020-01",
  "2020-02",
  "2020-03",
  "2020-04",
  "2020-05",
  "2020-06",
  "2020-07",
  "2020-08",
  "2020-09",
  "2021-10",
  "2021-11",
  "2021-12",
  "2021-01",
  "2021-02",
  "2021-03",
  "2021-04",
  "2021-05",
  "2021-06",
  "2021-07",
  "2021-08",
  "2021-09",
  "2022-10",
  "2022-11",
  "2022-12",
  "2022-01",
  "2022-02",
  "2022-03",
  "2022-04",
  "2022-05",
  "2022-06",
  "2022-07",
  "2022-08",
  "2022-09",
  "2023-01",
  "2023-02",
  "2023-03",
  "2023-04",
  "2023-05",
  "2023-06",
  "2023-07"
)

$values = @(
  53.9,
  54.6,
  54.6,
  54,
  54,
  54.7,
  53.9,
  54.1,
  54.7,
  54.1,
  53.7,
  55.1,
  53.1,
  52.8,
  52.6,
  54.6,
  52.9,
  54,
  54.1,
  54.6,
  54.4,
  53.6,
  53.8,
  54.1,
  52,
  53.7,
  53.4,
  53.2,
  52.4,
  54,
  53.4,
  53.3,
  53,
  53.1,
  53,
  53.1,
  55.3,
  55.7,
  55.1,
  53,
  28.9,
  53,
  53.4,
  53.4,
  54.2,
  54.1,
  54.5,
  55.1,
  50.8,
  52.2,
  52.2,
  52.8,
  51.6,
  55.3,
  53.8,
  54.2,
  52.9,
  52.4,
  48.9,
  51.7,
  49,
  47.1,
  42.6,
  51,
  51.2,
  48.8,
  42.7,
  48.4,
  54.1,
  52.5,
  51.7,
  50.9,
  52.9,
  56.4,
  57,
  54.4,
  52.9,
  52.3,
  51.1
)

for ($i = 0; $i -lt $dates.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $dates[$i]
  $ws.Cells.Item($row, 2).Value = $values[$i]
}